$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk-edit the "Weekly Groceries" transaction row (row 3) ---
# Date column: typo'd into a bogus value
$ws.Range("A3").Value = "20ddd26-01-12"
# Category column: typo'd
$ws.Range("C3").Value = "Groceriffes"
# Amount column: overwritten with stray text (no longer numeric)
$ws.Range("D3").Value = "as"

# --- Cosmetic/dashboard fixes ---
# Widen the Date and Description columns so the edited values are readable
$ws.Columns.Item(1).ColumnWidth = 28.625
$ws.Columns.Item(2).ColumnWidth = 21.125

# Leave the cursor on the cell that was being edited
$ws.Range("C3").Select()
